$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04271373187048221887
$ws.Range("C2").Value = 0.00006240767534437808
$ws.Range("D2").Value = 0.14942197473980470268
$ws.Range("E2").Value = 0.49423653606076972666
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.68643465034640094302
